$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 62. This shifts existing rows 62..145 down to 63..146,
# matching the target diff (dimension grows from A1:R145 to A1:R146).
$ws.Rows.Item(62).Insert()

# Populate the newly inserted row 62 with its values. Columns A, B, C, E, F, G, H, I,
# J, N, O, Q, R keep the same values the old row 62 had (now shifted to row 63),
# while D, K, L, M, P get the new figures from the update.
$ws.Cells.Item(62, 1).Value = 11
$ws.Cells.Item(62, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(62, 3).Value = "Bíobío"
$ws.Cells.Item(62, 4).Value = 44791
$ws.Cells.Item(62, 5).Value = 8
$ws.Cells.Item(62, 6).Value = 100112032
$ws.Cells.Item(62, 7).Value = "Zapallo italiano"
$ws.Cells.Item(62, 8).Value = "Sin especificar"
$ws.Cells.Item(62, 9).Value = "Primera"
$ws.Cells.Item(62, 10).Value = 220
$ws.Cells.Item(62, 11).Value = 21000
$ws.Cells.Item(62, 12).Value = 22000
$ws.Cells.Item(62, 13).Value = 21455
$ws.Cells.Item(62, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(62, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(62, 16).Value = 358
$ws.Cells.Item(62, 17).Value = 60
$ws.Cells.Item(62, 18).Value = "Hortaliza"

# Make sure the date cell keeps the date/time number format used by the rest of column D.
$ws.Cells.Item(62, 4).NumberFormat = $ws.Cells.Item(63, 4).NumberFormat
